$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 date_in to a real date value (was a text string "2/3/2020")
$ws.Range("D2").Value = 43892

# Update row 3 data: stock_code, no_of_shares, cost_per_share, date_in
$ws.Range("A3").Value = 1918
$ws.Range("B3").Value = 6000
$ws.Range("C3").Value = 32.65
$ws.Range("D3").Value = 44075

# Add new row 4
$ws.Range("A4").Value = 175
$ws.Range("B4").Value = 11000
$ws.Range("C4").Value = 16.18
$ws.Range("D4").Value = 44082

# Add new row 5
$ws.Range("A5").Value = 3883
$ws.Range("B5").Value = 11000
$ws.Range("C5").Value = 8.52
$ws.Range("D5").Value = 44082

# Apply the built-in short date number format (numFmtId 14) to D2, then
# replicate that exact same style onto D3:D5 via a formats-only paste so
# every cell in the date_in column shares a single cellXfs entry instead
# of each NumberFormat assignment minting its own style record.
$ws.Range("D2").NumberFormat = "mm-dd-yy"
$ws.Range("D2").Copy()
$ws.Range("D3:D5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D10").Select()
